$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "failed?" header to "status"
$ws.Range("E1").Value = "status"

# Add a new "error" header in column H
$ws.Range("H1").Value = "error"

# Add new row 18 of data (DEG run that errored out for SCTv2 corrected BL_A)
$ws.Range("A18").Value = "Neurolucida results"
$ws.Range("B18").Value = "2022-06-13 14-02-37"
$ws.Range("C18").Value = "DEG"
$ws.Range("D18").Value = "SCTv2 corrected BL_A"
$ws.Range("E18").Value = "error"
$ws.Range("F18").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G18").Value = "annotation, pseudotime"
$ws.Range("H18").Value = "Error in ValidateCellGroups(object = object, cells.1 = cells.1, cells.2 = cells.2,  : Cell group 2 is empty - no cells with identity class  Calls: sourceWithProgress ... FindMarkers -> FindMarkers.default -> ValidateCellGroups3"

# Add row 19 with a couple of spaces in F (leftover note)
$ws.Range("F19").Value = "  "

# Update the view: scroll so column B is the left-most visible column,
# and move/select the active cell to F22
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F22").Select()
